# === Week 3-1 introduction to css utility classes ===
# Regroups the CSS cheatsheet rows into three labeled sections
# - "стилове за измерения" (dimension utility classes: display, height, width,
#   border(-*), margin(-*), padding(-*))
# - "стилове за текст и шрифт" (text/font utility classes: text-decoration,
#   font-size, font-weight)
# - "стилове за цветове" (color utility classes: background-color, color)
# replacing the previous "цветове" / "измерения" / "текст и шрифт стилове" layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash the 6 distinct cell formats already present on the sheet into a
# scratch column (Z, well outside the table) before wiping the table, so they
# can be re-applied (format-only paste) once the content has been rebuilt in
# its new order/positions. (Scoping the later Clear() to A1:D30 keeps this
# scratch column intact.)
$ws.Range("A2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("A13").Copy()
$ws.Range("Z3").PasteSpecial(-4122)
$ws.Range("A23").Copy()
$ws.Range("Z4").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("Z5").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("Z6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Wipe the existing table content, formatting and merged regions.
$ws.Range("A1:D30").UnMerge()
$ws.Range("A1:D30").Clear()

# --- Re-enter every cell value in its new row position.
$ws.Range("A1").Value = "стилове за измерения"
$ws.Range("A2").Value = "display"
$ws.Range("B2").Value = "DIMENTION VALUE"
$ws.Range("A3").Value = "height"
$ws.Range("B3").Value = "NUMERIC + SIZE"
$ws.Range("A4").Value = "width"
$ws.Range("B4").Value = "NUMERIC + SIZE"
$ws.Range("A5").Value = "border"
$ws.Range("B5").Value = "TYPE COLOR NUMERIC + SIZE"
$ws.Range("A6").Value = "border-left"
$ws.Range("B6").Value = "TYPE COLOR NUMERIC + SIZE"
$ws.Range("A7").Value = "border-right"
$ws.Range("B7").Value = "TYPE COLOR NUMERIC + SIZE"
$ws.Range("A8").Value = "border-bottom"
$ws.Range("B8").Value = "TYPE COLOR NUMERIC + SIZE"
$ws.Range("A9").Value = "border-top"
$ws.Range("B9").Value = "TYPE COLOR NUMERIC + SIZE"
$ws.Range("A10").Value = "margin"
$ws.Range("B10").Value = "NUMERIC + SIZE"
$ws.Range("A11").Value = "margin-left"
$ws.Range("B11").Value = "NUMERIC + SIZE"
$ws.Range("A12").Value = "margin-right"
$ws.Range("B12").Value = "NUMERIC + SIZE"
$ws.Range("A13").Value = "margin-bottom"
$ws.Range("B13").Value = "NUMERIC + SIZE"
$ws.Range("A14").Value = "margin-top"
$ws.Range("B14").Value = "NUMERIC + SIZE"
$ws.Range("A15").Value = "padding"
$ws.Range("B15").Value = "NUMERIC + SIZE"
$ws.Range("A16").Value = "padding-left"
$ws.Range("B16").Value = "NUMERIC + SIZE"
$ws.Range("A17").Value = "padding-right"
$ws.Range("B17").Value = "NUMERIC + SIZE"
$ws.Range("A18").Value = "padding-bottom"
$ws.Range("B18").Value = "NUMERIC + SIZE"
$ws.Range("A19").Value = "padding-top"
$ws.Range("B19").Value = "NUMERIC + SIZE"
$ws.Range("A22").Value = "стилове за текст и шрифт"
$ws.Range("A23").Value = "text-decoration"
$ws.Range("B23").Value = "NONE"
$ws.Range("A24").Value = "font-size"
$ws.Range("B24").Value = "NUMERIC + SIZE"
$ws.Range("A25").Value = "font-weight"
$ws.Range("A27").Value = "стилове за цветове"
$ws.Range("A28").Value = "background-color"
$ws.Range("B28").Value = "COLOR"
$ws.Range("A29").Value = "color"
$ws.Range("B29").Value = "COLOR"

# --- Re-apply the stashed formats (format-only paste) to the cells that need them.
$ws.Range("Z5").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("Z5").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("Z4").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("Z6").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("Z6").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("Z5").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("Z5").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("Z5").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("Z5").Copy()
$ws.Range("B27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Clear the scratch stamp column now that every format has been transferred.
$ws.Range("Z1:Z6").Clear()

# --- Recreate the merged header / spacer cells.
$ws.Range("A1:B1").Merge()
$ws.Range("C5:D5").Merge()
$ws.Range("A22:B22").Merge()
$ws.Range("A27:B27").Merge()

# --- Restore the active selection cell.
$ws.Range("C10").Select()
